# Update results values (row 2) on each year sheet with the latest server
# results. Columns C, D, F, H, J, K remain 0 (unchanged).

$wb = $excel.ActiveWorkbook

$data = @{
    "2025" = @{
        A = 0
        B = 229.9102676639999
        E = 24938.4345679256
        G = 8095.925712661834
        I = 13878.11577882601
        L = 44963.47497589202
        M = 10112.3794939365
        N = 6479.486808609299
        O = 6069.678009745778
    }
    "2030" = @{
        A = 300.9804799738145
        B = 2946.312200738474
        E = 39741.37361765869
        G = 8095.925712661834
        I = 30652.64314713404
        L = 48981.18077820299
        M = 15653.53789831508
        N = 8247.748710722188
        O = 7089.117894102894
    }
    "2035" = @{
        A = 2135.626368860798
        B = 4546.031766517355
        E = 50486.27234148944
        G = 8095.925712661834
        I = 45851.45352390232
        L = 48981.18077820299
        M = 21186.37833031913
        N = 12213.55262468385
        O = 11755.72356163593
    }
    "2040" = @{
        A = 2135.626368860798
        B = 4546.031766517355
        E = 50486.27234148944
        G = 8095.925712661834
        I = 45851.45352390232
        L = 48981.18077820299
        M = 21186.37833031913
        N = 12213.55262468385
        O = 11755.72356163593
    }
    "2045" = @{
        A = 2135.626368860798
        B = 4546.031766517355
        E = 50486.27234148944
        G = 8095.925712661834
        I = 45851.45352390232
        L = 48981.18077820299
        M = 21186.37833031913
        N = 12213.55262468385
        O = 11755.72356163593
    }
    "2050" = @{
        A = 2135.626368860798
        B = 4546.031766517355
        E = 50486.27234148944
        G = 8095.925712661834
        I = 45851.45352390232
        L = 48981.18077820299
        M = 21186.37833031913
        N = 12213.55262468385
        O = 11755.72356163593
    }
}

foreach ($sheetName in $data.Keys) {
    $ws = $wb.Worksheets.Item([string]$sheetName)
    $rowValues = $data[$sheetName]
    foreach ($col in $rowValues.Keys) {
        $ws.Range("$col" + "2").Value = $rowValues[$col]
    }
}
